# Update Work Week and Social Spending
# --------------------------------------------------------------
# Refreshes the GDP per Capita series for Cameroon (Data sheet):
#   - revises the existing 1950-2010 figures to the latest release
#   - appends the newly published 2011-2016 figures
# The underlying source stores each "Data" figure as text (it is
# exported verbatim from the Clio-Infra dataset), so values are
# entered with a leading apostrophe to keep them as text instead of
# being auto-converted to numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Revise GDP per Capita values for existing years (1950-2010) ---
$ws.Range("E2").Value = "'1070"
$ws.Range("E3").Value = "'1095"
$ws.Range("E4").Value = "'1122"
$ws.Range("E5").Value = "'1148"
$ws.Range("E6").Value = "'1175"
$ws.Range("E7").Value = "'1202"
$ws.Range("E8").Value = "'1229"
$ws.Range("E9").Value = "'1256"
$ws.Range("E10").Value = "'1283"
$ws.Range("E11").Value = "'1310"
$ws.Range("E12").Value = "'1326"
$ws.Range("E13").Value = "'1321"
$ws.Range("E14").Value = "'1339"
$ws.Range("E15").Value = "'1366"
$ws.Range("E16").Value = "'1390"
$ws.Range("E17").Value = "'1393"
$ws.Range("E18").Value = "'1431"
$ws.Range("E19").Value = "'1443"
$ws.Range("E20").Value = "'1508"
$ws.Range("E21").Value = "'1551"
$ws.Range("E22").Value = "'1565"
$ws.Range("E23").Value = "'1578"
$ws.Range("E24").Value = "'1612"
$ws.Range("E25").Value = "'1599"
$ws.Range("E26").Value = "'1632"
$ws.Range("E27").Value = "'1677"
$ws.Range("E28").Value = "'1664"
$ws.Range("E29").Value = "'1706"
$ws.Range("E30").Value = "'1744"
$ws.Range("E31").Value = "'1785"
$ws.Range("E32").Value = "'1900"
$ws.Range("E33").Value = "'2153"
$ws.Range("E34").Value = "'2257"
$ws.Range("E35").Value = "'2345"
$ws.Range("E36").Value = "'2450"
$ws.Range("E37").Value = "'2585"
$ws.Range("E38").Value = "'2683"
$ws.Range("E39").Value = "'2471"
$ws.Range("E40").Value = "'2287"
$ws.Range("E41").Value = "'2021"
$ws.Range("E42").Value = "'1930"
$ws.Range("E43").Value = "'1833.35445177115"
$ws.Range("E44").Value = "'1754.2467380351"
$ws.Range("E45").Value = "'1677.14002614989"
$ws.Range("E46").Value = "'1615.59727367935"
$ws.Range("E47").Value = "'1648.64155827783"
$ws.Range("E48").Value = "'1710.10597530254"
$ws.Range("E49").Value = "'1780.49669668027"
$ws.Range("E50").Value = "'1847.2189118451"
$ws.Range("E51").Value = "'1902.26600633579"
$ws.Range("E52").Value = "'1961.08786225271"
$ws.Range("E53").Value = "'2028.54738427422"
$ws.Range("E54").Value = "'2084.9206379621"
$ws.Range("E55").Value = "'2143.98004322433"
$ws.Range("E56").Value = "'2198.23602014264"
$ws.Range("E57").Value = "'2222.91085339381"
$ws.Range("E58").Value = "'2269.5179519697"
$ws.Range("E59").Value = "'2315.60687676106"
$ws.Range("E60").Value = "'2352.79945998997"
$ws.Range("E61").Value = "'2370.13051001594"
$ws.Range("E62").Value = "'2420.1328496142"

# --- Append new rows for years 2011-2016 ---
$ws.Range("A63").Value = 120
$ws.Range("B63").Value = "Cameroon"
$ws.Range("C63").Value = "GDP per Capita"
$ws.Range("D63").Value = 2011
$ws.Range("E63").Value = "'2494"

$ws.Range("A64").Value = 120
$ws.Range("B64").Value = "Cameroon"
$ws.Range("C64").Value = "GDP per Capita"
$ws.Range("D64").Value = 2012
$ws.Range("E64").Value = "'2541"

$ws.Range("A65").Value = 120
$ws.Range("B65").Value = "Cameroon"
$ws.Range("C65").Value = "GDP per Capita"
$ws.Range("D65").Value = 2013
$ws.Range("E65").Value = "'2612"

$ws.Range("A66").Value = 120
$ws.Range("B66").Value = "Cameroon"
$ws.Range("C66").Value = "GDP per Capita"
$ws.Range("D66").Value = 2014
$ws.Range("E66").Value = "'2696"

$ws.Range("A67").Value = 120
$ws.Range("B67").Value = "Cameroon"
$ws.Range("C67").Value = "GDP per Capita"
$ws.Range("D67").Value = 2015
$ws.Range("E67").Value = "'2778"

$ws.Range("A68").Value = 120
$ws.Range("B68").Value = "Cameroon"
$ws.Range("C68").Value = "GDP per Capita"
$ws.Range("D68").Value = 2016
$ws.Range("E68").Value = "'2828"

